# Natmi following Dr Hou advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ "E"=3; "G"=15.140316; "H"=45.420948; "I"=0.3229157245229468; "J"=0.3229157245229468; "K"=3; "M"=0.7567149999999999; "N"=2.270145; "O"=0.07973436174489927; "P"=0.07973436174489928; "Q"=11.45690422194; "R"=103.11213799746; "S"=0.02574747919222888; "T"=0.02574747919222888 }
    3 = @{ "E"=3; "G"=15.140316; "H"=45.420948; "I"=0.3229157245229468; "J"=0.3229157245229468; "K"=3; "M"=3.474244666666667; "N"=10.422734; "O"=0.3660779567502785; "P"=0.3660779567502786; "Q"=52.601162114648; "R"=473.410459031832; "S"=0.1182123286358962; "T"=0.1182123286358962 }
    4 = @{ "E"=3; "G"=15.140316; "H"=45.420948; "I"=0.3229157245229468; "J"=0.3229157245229468; "K"=3; "M"=2.264125333333333; "N"=6.792376; "O"=0.2385687985090697; "P"=0.2385687985090697; "Q"=34.279573010272; "R"=308.516157092448; "S"=0.07703761641912514; "T"=0.07703761641912514 }
    5 = @{ "E"=3; "G"=15.140316; "H"=45.420948; "I"=0.3229157245229468; "J"=0.3229157245229468; "K"=3; "M"=2.995365333333333; "N"=8.986096; "O"=0.3156188829957525; "P"=0.3156188829957525; "Q"=45.350777682112; "R"=408.156999139008; "S"=0.1019183002756966; "T"=0.1019183002756966 }
    6 = @{ "E"=3; "G"=18.94069966666667; "H"=56.822099; "I"=0.4039710766824948; "J"=0.4039710766824948; "K"=3; "M"=0.7567149999999999; "N"=2.270145; "O"=0.07973436174489927; "P"=0.07973436174489928; "Q"=14.33271154826167; "R"=128.994403934355; "S"=0.03221037596267848; "T"=0.03221037596267849 }
    7 = @{ "E"=3; "G"=18.94069966666667; "H"=56.822099; "I"=0.4039710766824948; "J"=0.4039710766824948; "K"=3; "M"=3.474244666666667; "N"=10.422734; "O"=0.3660779567502785; "P"=0.3660779567502786; "Q"=65.80462479985178; "R"=592.241623198666; "S"=0.1478849063381378; "T"=0.1478849063381378 }
    8 = @{ "E"=3; "G"=18.94069966666667; "H"=56.822099; "I"=0.4039710766824948; "J"=0.4039710766824948; "K"=3; "M"=2.264125333333333; "N"=6.792376; "O"=0.2385687985090697; "P"=0.2385687985090697; "Q"=42.88411794635822; "R"=385.957061517224; "S"=0.09637489439655805; "T"=0.09637489439655805 }
    9 = @{ "E"=3; "G"=18.94069966666667; "H"=56.822099; "I"=0.4039710766824948; "J"=0.4039710766824948; "K"=3; "M"=2.995365333333333; "N"=8.986096; "O"=0.3156188829957525; "P"=0.3156188829957525; "Q"=56.73431517061156; "R"=510.608836535504; "S"=0.1275008999851205; "T"=0.1275008999851205 }
    10 = @{ "E"=3; "G"=3.221232; "H"=9.663696; "I"=0.06870308817441464; "J"=0.06870308817441464; "K"=3; "M"=0.7567149999999999; "N"=2.270145; "O"=0.07973436174489927; "P"=0.07973436174489928; "Q"=2.43755457288; "R"=21.93799115592; "S"=0.005477996885490489; "T"=0.00547799688549049 }
    11 = @{ "E"=3; "G"=3.221232; "H"=9.663696; "I"=0.06870308817441464; "J"=0.06870308817441464; "K"=3; "M"=3.474244666666667; "N"=10.422734; "O"=0.3660779567502785; "P"=0.3660779567502786; "Q"=11.191348096096; "R"=100.722132864864; "S"=0.02515068614132394; "T"=0.02515068614132394 }
    12 = @{ "E"=3; "G"=3.221232; "H"=9.663696; "I"=0.06870308817441464; "J"=0.06870308817441464; "K"=3; "M"=2.264125333333333; "N"=6.792376; "O"=0.2385687985090697; "P"=0.2385687985090697; "Q"=7.293272975744; "R"=65.639456781696; "S"=0.01639041319963278; "T"=0.01639041319963278 }
    13 = @{ "E"=3; "G"=3.221232; "H"=9.663696; "I"=0.06870308817441464; "J"=0.06870308817441464; "K"=3; "M"=2.995365333333333; "N"=8.986096; "O"=0.3156188829957525; "P"=0.3156188829957525; "Q"=9.648766663424; "R"=86.838899970816; "S"=0.02168399194796744; "T"=0.02168399194796744 }
    14 = @{ "E"=3; "G"=9.584029000000001; "H"=28.752087; "I"=0.2044101106201438; "J"=0.2044101106201438; "K"=3; "M"=0.7567149999999999; "N"=2.270145; "O"=0.07973436174489927; "P"=0.07973436174489928; "Q"=7.252378504735; "R"=65.271406542615; "S"=0.01629850970450142; "T"=0.01629850970450143 }
    15 = @{ "E"=3; "G"=9.584029000000001; "H"=28.752087; "I"=0.2044101106201438; "J"=0.2044101106201438; "K"=3; "M"=3.474244666666667; "N"=10.422734; "O"=0.3660779567502785; "P"=0.3660779567502786; "Q"=33.29726163842867; "R"=299.675354745858; "S"=0.07483003563492066; "T"=0.07483003563492067 }
    16 = @{ "E"=3; "G"=9.584029000000001; "H"=28.752087; "I"=0.2044101106201438; "J"=0.2044101106201438; "K"=3; "M"=2.264125333333333; "N"=6.792376; "O"=0.2385687985090697; "P"=0.2385687985090697; "Q"=21.69944285430133; "R"=195.294985688712; "S"=0.04876587449375373; "T"=0.04876587449375374 }
    17 = @{ "E"=3; "G"=9.584029000000001; "H"=28.752087; "I"=0.2044101106201438; "J"=0.2044101106201438; "K"=3; "M"=2.995365333333333; "N"=8.986096; "O"=0.3156188829957525; "P"=0.3156188829957525; "Q"=28.70766822026134; "R"=258.369013982352; "S"=0.06451569078696799; "T"=0.06451569078696799 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
